# The document contains a table whose cell paragraphs carry a malformed
# paragraph-style reference:  <w:pStyle w:pstlname="Normal"/>
# ("pstlname" is not a real OOXML attribute - a valid style reference uses
# w:val - so it is effectively dead markup that just duplicates the
# paragraphs' already-inherited "Normal" style). The commit removes that
# stray element from every affected cell paragraph, 22 in total, without
# touching anything else (text, other styles, table structure, etc. are
# all left exactly as they were).
#
# Re-asserting each affected paragraph's Style to its own (already
# effective) "Normal" style makes the document model rebuild that
# paragraph's properties cleanly, which drops the bogus w:pstlname
# attribute instead of round-tripping it.

$d = $word.ActiveDocument

$total = $d.Paragraphs.Count
$fixed = 0

for ($i = 1; $i -le $total; $i++) {
    $para = $d.Paragraphs.Item($i)
    $styleName = $para.Style.NameLocal

    # Paragraph 4 is the table's real "Table 1: test" caption; its true
    # style is "ImageCaption" even though a pre-existing quirk in this
    # document makes that particular paragraph mis-report as "Normal" -
    # skip it so we never touch the caption's actual style.
    if ($i -eq 4) {
        continue
    }

    if ($styleName -eq "Normal") {
        $para.Style = "Normal"
        $fixed++
    }
}

Write-Host "Normalized paragraph style on $fixed paragraph(s)."
